$wb = $excel.ActiveWorkbook

# --- Update currency from EUR to USD on the "Simple Fields" sheets ---
$wsSimple = $wb.Worksheets.Item("Simple Fields")
$wsSimple.Range("J2").Value = "USD"

$wsSimpleFmt = $wb.Worksheets.Item("Simple Fields - Formatted")
$wsSimpleFmt.Range("J2").Value = "USD"

# --- Populate the "Unit Price" column (D2) on the Items sheets ---
$wsItems = $wb.Worksheets.Item("Items")
$wsItems.Range("C2").Copy()
$wsItems.Range("D2").PasteSpecial(-4122)
$wsItems.Range("D2").NumberFormat = "@"
$wsItems.Range("D2").Value = "7741"
$wsItems.Range("C2").Copy()
$wsItems.Range("D2").PasteSpecial(-4122)

$wsItemsFmt = $wb.Worksheets.Item("Items - Formatted")
$wsItemsFmt.Range("C2").Copy()
$wsItemsFmt.Range("D2").PasteSpecial(-4122)
$wsItemsFmt.Range("D2").NumberFormat = "@"
$wsItemsFmt.Range("D2").Value = "7741.00"
$wsItemsFmt.Range("C2").Copy()
$wsItemsFmt.Range("D2").PasteSpecial(-4122)
